# Fruta / hortaliza, semanal
# Insert 3 new weekly records (Femacal de La Calera - Ciruela, variedad
# "Angeleno", unidad "$/bandeja 10 kilos granel") above the existing
# row 293, pushing the previously existing rows 293-321 down to 296-324.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three blank rows at the top of the block (row 293).
$ws.Rows.Item(293).Insert()
$ws.Rows.Item(293).Insert()
$ws.Rows.Item(293).Insert()

# --- New row 293 ---------------------------------------------------------
$ws.Cells.Item(293,1).Value  = 3
$ws.Cells.Item(293,2).Value  = "Femacal de La Calera"
$ws.Cells.Item(293,3).Value  = "Coquimbo"
$ws.Cells.Item(293,4).Value  = [datetime]"2023-03-28"
$ws.Cells.Item(293,5).Value  = 5
$ws.Cells.Item(293,6).Value  = "Fruta"
$ws.Cells.Item(293,7).Value  = 100103
$ws.Cells.Item(293,8).Value  = "Frutos de hueso (carozo)"
$ws.Cells.Item(293,9).Value  = 100103002
$ws.Cells.Item(293,10).Value = "Ciruela"
$ws.Cells.Item(293,11).Value = "Angeleno"
$ws.Cells.Item(293,12).Value = "Especial"
$ws.Cells.Item(293,13).Value = 48
$ws.Cells.Item(293,14).Value = 10000
$ws.Cells.Item(293,15).Value = 10000
$ws.Cells.Item(293,16).Value = 10000
$ws.Cells.Item(293,17).Value = "`$/bandeja 10 kilos granel"
$ws.Cells.Item(293,18).Value = "Región de O'Higgins"
$ws.Cells.Item(293,19).Value = 1000
$ws.Cells.Item(293,20).Value = 10

# --- New row 294 ---------------------------------------------------------
$ws.Cells.Item(294,1).Value  = 3
$ws.Cells.Item(294,2).Value  = "Femacal de La Calera"
$ws.Cells.Item(294,3).Value  = "Coquimbo"
$ws.Cells.Item(294,4).Value  = [datetime]"2023-03-28"
$ws.Cells.Item(294,5).Value  = 5
$ws.Cells.Item(294,6).Value  = "Fruta"
$ws.Cells.Item(294,7).Value  = 100103
$ws.Cells.Item(294,8).Value  = "Frutos de hueso (carozo)"
$ws.Cells.Item(294,9).Value  = 100103002
$ws.Cells.Item(294,10).Value = "Ciruela"
$ws.Cells.Item(294,11).Value = "Angeleno"
$ws.Cells.Item(294,12).Value = "Primera"
$ws.Cells.Item(294,13).Value = 50
$ws.Cells.Item(294,14).Value = 8000
$ws.Cells.Item(294,15).Value = 8000
$ws.Cells.Item(294,16).Value = 8000
$ws.Cells.Item(294,17).Value = "`$/bandeja 10 kilos granel"
$ws.Cells.Item(294,18).Value = "Región de O'Higgins"
$ws.Cells.Item(294,19).Value = 800
$ws.Cells.Item(294,20).Value = 10

# --- New row 295 ---------------------------------------------------------
$ws.Cells.Item(295,1).Value  = 3
$ws.Cells.Item(295,2).Value  = "Femacal de La Calera"
$ws.Cells.Item(295,3).Value  = "Coquimbo"
$ws.Cells.Item(295,4).Value  = [datetime]"2023-03-28"
$ws.Cells.Item(295,5).Value  = 5
$ws.Cells.Item(295,6).Value  = "Fruta"
$ws.Cells.Item(295,7).Value  = 100103
$ws.Cells.Item(295,8).Value  = "Frutos de hueso (carozo)"
$ws.Cells.Item(295,9).Value  = 100103002
$ws.Cells.Item(295,10).Value = "Ciruela"
$ws.Cells.Item(295,11).Value = "Angeleno"
$ws.Cells.Item(295,12).Value = "Segunda"
$ws.Cells.Item(295,13).Value = 47
$ws.Cells.Item(295,14).Value = 6000
$ws.Cells.Item(295,15).Value = 6000
$ws.Cells.Item(295,16).Value = 6000
$ws.Cells.Item(295,17).Value = "`$/bandeja 10 kilos granel"
$ws.Cells.Item(295,18).Value = "Región de O'Higgins"
$ws.Cells.Item(295,19).Value = 600
$ws.Cells.Item(295,20).Value = 10
